$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week02 task for student PPF004 (row 24) is now complete:
#  - mark name + status cells with the "Completed" highlight style (copy the
#    formatting already used for other completed rows, e.g. row 2)
#  - clear the now-irrelevant "Pending Task" cell (D24)
#  - flip the status text to "Completed" and bump the streak to 1

$ws.Range("B2").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("D24").ClearContents()
$ws.Range("E24").Value = "Completed"
$ws.Range("F24").Value = 1

# Leave the selection where the user ended up after editing the row.
$ws.Range("E24").Select()
